# Updated cryptos list (Price and Volume(1h) columns) to match the
# scheduled GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the columns whose new "Price" values look like plain numbers to
# stay stored as text (matches the source data, which is text, not numeric).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.980.57"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "1.639.89"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  -0.66%  "
$ws.Range("D5").Value = "214.95"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").Value = "0.5099"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("D8").Value = "0.2582"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").Value = "0.06352"
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").Value = "19.78"
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("D11").Value = "0.07762"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").Value = "4.276"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("D13").Value = "1.636.65"
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("D14").Value = "0.5464"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").Value = "0.0₅7740"
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("D17").Value = "26.003.74"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").Value = "196.23"
$ws.Range("E19").Value = "  -1.51%  "
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").Value = "9.915"
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("D22").Value = "6.077"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("D24").Value = "1.896"
$ws.Range("E24").Value = "  +1.63%  "
$ws.Range("D25").Value = "142.97"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("D26").Value = "0.1229"
$ws.Range("E26").Value = "  +6.61%  "
$ws.Range("D27").Value = "6.861"
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").Value = "15.58"
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("D30").Value = "0.04867"
$ws.Range("E30").Value = "  -3.38%  "
$ws.Range("D31").Value = "3.276"
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("D32").Value = "3.215"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("D34").Value = "2.377"
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("D35").Value = "0.9143"
$ws.Range("E35").Value = "  +1.93%  "
$ws.Range("D36").Value = "2.563"
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").Value = "1.091.63"
$ws.Range("E38").Value = "  -4.26%  "
$ws.Range("D39").Value = "0.01568"
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D41").Value = "2.528"
$ws.Range("E41").Value = "  -1.34%  "
$ws.Range("D42").Value = "5.584"
$ws.Range("E42").Value = "  -1.63%  "
$ws.Range("D43").Value = "0.8052"
$ws.Range("E43").Value = "  -1.59%  "
$ws.Range("D44").Value = "99.19"
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("E45").Value = "  -2.49%  "
$ws.Range("D46").Value = "1.780.02"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").Value = "1.007"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "55.16"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").Value = "0.05213"
$ws.Range("E50").Value = "  +2.24%  "
$ws.Range("D51").Value = "7.497"
$ws.Range("E51").Value = "  +0.70%  "
